$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.346130013465881
$ws.Range("B1").Value = 1.478269577026367
$ws.Range("C1").Value = 3.953767538070679
$ws.Range("D1").Value = 3.208298444747925
$ws.Range("E1").Value = 1.079349517822266
